# Commit: "accept multilple filter keys"
# Insert a new (empty) column B on the "Desikan" and "Destrieux" sheets,
# shifting the existing "ABCD ROI name"/value columns one slot to the right
# so users can filter/match on more than one key column. The "ASEG" sheet
# is left structurally alone; only view/selection state changes.

$wb = $excel.ActiveWorkbook

# ---- Desikan (sheet1): insert blank column B ----------------------------
$ws1 = $wb.Worksheets.Item("Desikan")
$ws1.Columns.Item(2).Insert()
$ws1.Columns.Item(2).ColumnWidth = 18.43

# ---- Destrieux (sheet2): insert blank column B ---------------------------
$ws2 = $wb.Worksheets.Item("Destrieux")
$ws2.Columns.Item(2).Insert()
$ws2.Columns.Item(2).ColumnWidth = 18.43

# ---- ASEG (sheet3): unchanged structurally -------------------------------
$ws3 = $wb.Worksheets.Item("ASEG")

# ---- View / selection state ----------------------------------------------
# Select/activate in the same order Excel would have recorded, finishing on
# "Destrieux" so it ends up as the tabSelected / active sheet.
$ws1.Activate()
$ws1.Range("B9").Select()

$ws2.Activate()
$ws2.Range("C149").Select()

$ws3.Activate()
$ws3.Range("B8").Select()

$ws2.Activate()
